# 100. Same Tree (Trees-easy)
#
# Mirrors the existing "Trees ( Basics )" block (section-header row 82 +
# data rows 84/85) with a new "Trees ( Advance )" block: section header in
# row 87, and a new question entry (Q No. 100, "Same Tree") in row 89.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 87: new section header, styled like the other section headers
# (e.g. B82 "Trees ( Basics )" / B78 / B62 / ...). Copy the format first
# (keeps the existing shared cell style instead of minting a new one),
# then write the text.
$ws.Range("B82").Copy() | Out-Null
$ws.Range("B87").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B87").Value = "Trees ( Advance )"

# --- Row 89: new question row, styled like the other question rows
# (e.g. row 85 "Maximum Depth of Binary Tree" / row 84 "Invert Binary Tree").
$ws.Range("A89").Value = 100

$ws.Range("B85").Copy() | Out-Null
$ws.Range("B89").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B89").Value = "Same Tree"

$ws.Range("C89").Value = "Easy"
$ws.Range("D89").Value = "DFS,recursion"

# Date solved: 2025-05-07 -> serial 45784 (same "General" display the
# other Date-solved cells use).
$ws.Range("E89").Value = 45784

$excel.CutCopyMode = $false

# Update the view to match the new bottom of the sheet, like the original
# author's selection move after adding the rows.
$ws.Range("B90").Select()
